$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (existing H:Q shift right to I:R),
# mirroring the author's "add a supplier-name column" edit. Insert() also
# carries the per-row cell styles for the shifted cells along correctly.
$ws.Columns("H").Insert()

# Match the new column's width to its left neighbour (G) so the pair reads
# as one formatted block, then give the new column its header text.
$ws.Columns("H").ColumnWidth = $ws.Columns("G").ColumnWidth
$ws.Range("H1").Value = "Tên NCC"

# The autofilter range is stale after the insert (still only reaches the
# old last column); turn it off and re-apply across the full new range.
$ws.AutoFilterMode = $false
$ws.Range("A1:Q2").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# autofilter's new range.
$sheetName = $ws.Name
foreach ($n in $wb.Names) {
    $n.RefersTo = "='" + $sheetName + "'!`$A`$1:`$Q`$2"
}

# Move the selection, matching the edited workbook's saved cursor position.
$ws.Range("H4").Select()
